# Applies the OOXML diff:
#   - Several paragraphs had their text split across multiple <w:r> runs
#     with <w:proofErr w:type="spellStart"/>...<w:proofErr w:type="spellEnd"/>
#     markers wrapping the "misspelled" word (a normal side-effect of
#     Word's spell checker). The edit removes those spell-check markers
#     and the run splits they caused, leaving a single merged run with
#     the identical visible text.
#   - One paragraph ("Pobranie obecnie zamówionych kursów") gets new
#     trailing text " (kierowca)" appended.
#
# Word's Find/Replace, when the search hits a range that spans several
# runs (possibly interspersed with <w:proofErr/> markers), rewrites the
# matched range as a single run and drops the now-orphaned proofErr
# markers - exactly the behaviour the diff shows, so a simple
# self-replace ("find the visible text, replace it with itself") is
# enough to normalize each of those paragraphs.

$d = $word.ActiveDocument

function Merge-Runs($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. "UITMBER (" + "WSIiZ" + "/UITM UBER)" -> single run
Merge-Runs "UITMBER `(WSIiZ/UITM UBER`)" "UITMBER (WSIiZ/UITM UBER)"

# 2. "Data " + "Driven" + " Development" -> single run
Merge-Runs "Data Driven Development" "Data Driven Development"

# 3. Xamarin.Forms run: drop the spellStart/spellEnd wrap (text unchanged)
Merge-Runs "Xamarin.Forms" "Xamarin.Forms"

# 4. "AppleId" + "," -> single run
Merge-Runs "AppleId," "AppleId,"

# 5. "Typ auta (Uber, " + "UberX" + ", 7-osobowe auto)" -> single run
Merge-Runs "Typ auta `(Uber, UberX, 7-osobowe auto`)" "Typ auta (Uber, UberX, 7-osobowe auto)"

# 6. "Historia ... (cena, " + "kierwca" + ", czas, trasa" -> single run
#    (the trailing ", klient" and ")" runs are untouched / stay separate)
Merge-Runs "Historia przejazdu ze szczegółami `(cena, kierwca, czas, trasa" `
           "Historia przejazdu ze szczegółami (cena, kierwca, czas, trasa"

# 7. " spięcia z " + "Spotify" -> single run ("Możliwość" run stays separate)
Merge-Runs " spięcia z Spotify" " spięcia z Spotify"

# 8. " (" + "zdjecie" + ", auto, rejestracja)" -> single run
#    ("Informacje o kierowcy" run stays separate)
Merge-Runs " `(zdjecie, auto, rejestracja`)" " (zdjecie, auto, rejestracja)"

# 9. ".Net " + "Core" + " " + "Api" -> single run
Merge-Runs ".Net Core Api" ".Net Core Api"

# 10. Append " (kierowca)" right after "Pobranie obecnie zamówionych kursów"
$tail = $d.Content
$found = $tail.Find.Execute("Pobranie obecnie zamówionych kursów", $true, $false, `
                             $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $tail.Collapse(0)
    $tail.InsertAfter(" (kierowca)")
}
